$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (previously "Strike#") values. Regenerate them for rows 2-7.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 1
